$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "title"
$ws.Range("C3").Value = "input"
$ws.Range("C4").Value = "input"
$ws.Range("C5").Value = "input"
$ws.Range("C6").Value = "description"

$ws.Range("C2").Select()
